$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "80239121X"
$ws.Range("C17").Value = "Gallego Doncel, Alejandro"
$ws.Range("D17").Value = "UO285577@uniovi.es"
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:UO285577@uniovi.es", "", "", "UO285577@uniovi.es")
$ws.Range("D17").Font.Underline = $true
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = "No"
$ws.Range("H17").Value = "No"
$ws.Range("I17").Value = "Clases Expositivas-A"
$ws.Range("J17").Value = "Practicas de Aula/Semin-01"
$ws.Range("K17").Value = "Prácticas de Laboratorio-01"
$ws.Range("L17").Value = "Tutorías Grupales-02"

$ws.Range("L18").Select() | Out-Null
